# Table_Q38: fix column ordering for the multi-answer (Yes/No) table.
#
# The table has 5 grid columns:
#   1: Question   2-3: Ecology (No | Yes)   4-5: Social (No | Yes)
#
# In the header sub-row the "No"/"Yes" labels were generated in the
# wrong order (No, Yes) and likewise the counts in the data row were
# generated in the wrong order (26, 3 / 12, 1). The fix swaps each
# pair of columns so the labels read (Yes, No) and the counts line up
# with them.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Get-CellText($table, $row, $col) {
    return $table.Cell($row, $col).Range.Text.TrimEnd([char]13, [char]7)
}

# --- Row 2: the "No" / "Yes" sub-header ------------------------------
$row = 2
$c2 = Get-CellText $t $row 2
$c3 = Get-CellText $t $row 3
$c4 = Get-CellText $t $row 4
$c5 = Get-CellText $t $row 5

$t.Cell($row, 2).Range.Text = $c3
$t.Cell($row, 3).Range.Text = $c2
$t.Cell($row, 4).Range.Text = $c5
$t.Cell($row, 5).Range.Text = $c4

# --- Row 3: the answer counts -----------------------------------------
$row = 3
$c2 = Get-CellText $t $row 2
$c3 = Get-CellText $t $row 3
$c4 = Get-CellText $t $row 4
$c5 = Get-CellText $t $row 5

$t.Cell($row, 2).Range.Text = $c3
$t.Cell($row, 3).Range.Text = $c2
$t.Cell($row, 4).Range.Text = $c5
$t.Cell($row, 5).Range.Text = $c4
